$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the existing data (row 79)
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# New quote date: 2025-11-23 (serial 45984)
$ws.Cells.Item($newRow, 1).Value = 45984
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($newRow, 2).Value = "21,4507"
$ws.Cells.Item($newRow, 3).Value = "15,7543"
$ws.Cells.Item($newRow, 4).Value = "15,2104"
$ws.Cells.Item($newRow, 5).Value = "15,2104"
